# Add a new "PRESUPUESTO" column (G) to the "VENTA MENSUAL" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# Copy the formatting (styles/number formats/borders) from column F into the
# new column G, then overwrite the contents with the new header and values.
$ws.Range("F1:F19").Copy()
$ws.Range("G1:G19").PasteSpecial(-4122)

# Column width for the new column G (target stored width of 17 characters;
# Excel's ColumnWidth property is offset from the stored XML width by the
# standard ~0.83 character padding used by this workbook's column widths)
$ws.Columns.Item(7).ColumnWidth = 16.17

# Header in G1
$ws.Range("G1").Value = "PRESUPUESTO"

# Fill G2:G19 with 0 (budget values, all zero)
for ($r = 2; $r -le 19; $r++) {
    $ws.Cells.Item($r, 7).Value = 0
}
